$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.500.44"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.732.21"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4893"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2675"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06220"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "1.732.84"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07065"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.69"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.650"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6099"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.43"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D17").Value = "26.492.41"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007181"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.94%  "
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("D21").Value = "1.957.30"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.785"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.260"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.47"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.776"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "108.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.403"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.975"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08045"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.697"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04583"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.617"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.008"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6379"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9001"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.021"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.004"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01509"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.29"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.452"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3896"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.950"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1185"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05385"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.799"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.250"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3413"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.19%  "
